# Daily update at 8 AM UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous "latest" row (26) loses its special date-only number format
# and reverts to the regular date+time format used by all earlier rows.
$ws.Range("A26").NumberFormat = $ws.Range("A25").NumberFormat

# Append the new day's results as row 27, giving it the "latest" date-only
# number format that row 26 used to have.
$ws.Range("A27").Value = 45767
$ws.Range("A27").NumberFormat = "YYYY-MM-DD"
$ws.Range("B27").Value = 110
$ws.Range("C27").Value = 111
$ws.Range("D27").Value = 108
